# Regenerate merged AHB files
# 1) Rename the header labels in row 1 from *_old / *_new to *_FV2210 / *_FV2304
# 2) Turn the data range into an Excel Table ("Table1")
# 3) Freeze the header row (row 1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2210","Segmentgruppe_FV2210","Segment_FV2210","Datenelement_FV2210","Segment ID_FV2210",
  "Code_FV2210","Qualifier_FV2210","Beschreibung_FV2210","Bedingungsausdruck_FV2210","Bedingung_FV2210",
  "diff",
  "Segmentname_FV2304","Segmentgruppe_FV2304","Segment_FV2304","Datenelement_FV2304","Segment ID_FV2304",
  "Code_FV2304","Qualifier_FV2304","Beschreibung_FV2304","Bedingungsausdruck_FV2304","Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert the used range into a native Excel table
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), [System.Reflection.Missing]::Value, 1)
$lo.Name = "Table1"

# Freeze the header row
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
